$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: replace numeric/date header cells with text labels ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "attendance time"
$ws.Range("C1").Value = "departure time"
$ws.Range("D1").Value = "attendance Date"

$a1 = $ws.Range("A1")
$a1.ClearFormats()
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108

$a1.Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2: fix the attendance-date (column D) value ---
$ws.Range("D2").Value = 44502

# --- Rows 3-10: new seed data (id, attendance time, departure time, attendance date) ---
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 44542.333333333336
$ws.Range("C3").Value = 44542.708333333336
$ws.Range("D3").Value = 44502

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 44542.333333333336
$ws.Range("C4").Value = 44542.666666666664
$ws.Range("D4").Value = 44502

$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 44542.416666666664
$ws.Range("C5").Value = 44542.708333333336
$ws.Range("D5").Value = 44502

$ws.Range("A6").Value = 6
$ws.Range("B6").Value = 44542.416666666664
$ws.Range("C6").Value = 44542.625
$ws.Range("D6").Value = 44502

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = 44542.416666666664
$ws.Range("C7").Value = 44542.666666666664
$ws.Range("D7").Value = 44502

$ws.Range("A8").Value = 8
$ws.Range("B8").Value = 44542.375
$ws.Range("C8").Value = 44542.666666666664
$ws.Range("D8").Value = 44502

$ws.Range("A9").Value = 9
$ws.Range("B9").Value = 44542.375
$ws.Range("C9").Value = 44542.708333333336
$ws.Range("D9").Value = 44502

$ws.Range("A10").Value = 10
$ws.Range("B10").Value = 44542.375
$ws.Range("C10").Value = 44542.625
$ws.Range("D10").Value = 44502
$ws.Range("B10:D10").NumberFormat = "m/d/yy h:mm"

# --- Selection matches the saved file: D10 active ---
[void]$ws.Range("D10").Select()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
